# Update each sheet's remn_amt (column B) table:
#  - fill in B100 and B101 (previously placeholder 0 values)
#  - append a new row 102 (date 45961, remn_amt 0)
#
# Values taken from the target diff, one triple (B100, B101, B102-date) per
# worksheet, in sheet/tab order.

$wb = $excel.ActiveWorkbook

$values = @(
    @{ B100 = 10319015; B101 = 10671656 },
    @{ B100 = 13415088; B101 = 13473681 },
    @{ B100 = 3808734;  B101 = 3596902  },
    @{ B100 = 1070636;  B101 = 1030573  },
    @{ B100 = 1901071;  B101 = 1819220  },
    @{ B100 = 1980666;  B101 = 1935706  },
    @{ B100 = 300917;   B101 = 303355   },
    @{ B100 = 322170;   B101 = 332237   }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $v = $values[$i]

    $ws.Cells.Item(100, 2).Value = $v.B100
    $ws.Cells.Item(101, 2).Value = $v.B101

    $ws.Cells.Item(102, 1).Value = 45961
    $ws.Cells.Item(102, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item(102, 2).Value = 0
}
